$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "CasesTab" row (row 2) query text used to return a `Cohort` column
# derived from `co.cohort_description`. That trailing RETURN clause /
# OPTIONAL MATCH usage was erroneous for this query and has been removed,
# matching the corrected query already used elsewhere in the workbook.
$fixedCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Bouvier des Flandres']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $fixedCasesQuery

# Match the author's final view state: selection on B2, scrolled back to
# the top of the sheet (no forced topLeftCell override).
$ws.Range("B2").Select() | Out-Null

# Row heights shrink slightly because row 2's text now has fewer wrapped
# lines (matching the unchanged FilesTab query height) while the sheet's
# general text-metrics also tightened a touch; reapply the observed
# auto-fit heights for all three data rows.
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 244.8
